$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-7 (columns A..T). Row 1 is the header and is unchanged.
# Two new rows (6 and 7) are appended, and rows 2-5 are updated in place so
# that the whole block reflects the latest weekly price records.
$data = @(
    @(11, "Vega Monumental Concepción", "Bíobío", 44545, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Castle Brite", "Primera", 100, 18000, 19000, 18500, "`$/caja 15 kilos", "Región de O'Higgins", 1233, 15),
    @(11, "Vega Monumental Concepción", "Bíobío", 44545, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Castle Brite", "Segunda", 50, 17000, 17000, 17000, "`$/caja 15 kilos", "Región de O'Higgins", 1133, 15),
    @(11, "Vega Monumental Concepción", "Bíobío", 44159, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Castle Brite", "Primera", 100, 14000, 15000, 14500, "`$/caja 15 kilos", "Región Metropolitana", 967, 15),
    @(11, "Vega Monumental Concepción", "Bíobío", 44189, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Dina", "Primera", 200, 15000, 16000, 15500, "`$/caja 15 kilos granel", "Región de O'Higgins", 1033, 15),
    @(11, "Vega Monumental Concepción", "Bíobío", 44189, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Dina", "Segunda", 100, 14000, 14000, 14000, "`$/caja 15 kilos granel", "Región de O'Higgins", 933, 15),
    @(11, "Vega Monumental Concepción", "Bíobío", 44187, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Dina", "Primera", 100, 15000, 16000, 15500, "`$/caja 18 kilos", "Región Metropolitana", 861, 18)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

# Rows 6 and 7 are brand new; give their "Fecha" (date) cells the same
# number format used by the rest of the column (column D).
$ws.Range("D6:D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
